$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
# -----------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Experience the Montreal circus-themed Cirque Du Soleil Kooza slot game for free. Read our review covering gameplay, graphics, theme, and pro and cons.</w:t></w:r></w:p>'
$metaPara.Range.InsertXML($metaXml)

# -----------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the
#    document, and replace the remaining italic paragraph's text with
#    the DALLE image prompt (keeping its italic run formatting).
# -----------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete()

$count = $d.Paragraphs.Count
$descPara = $d.Paragraphs.Item($count)
$descTextRange = $d.Range($descPara.Range.Start, $descPara.Range.End - 1)
$descTextRange.Text = "Prompt: DALLE, please create a feature image for Cirque du Soleil Kooza that captures the whimsical and colorful nature of the game and its circus theme. The image should be in a cartoon style and feature a happy Maya warrior with glasses. Make sure it is eye-catching and reflects the excitement and fun of this slot game."
